$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 127, pushing the existing rows 127-154 down to 130-157.
$ws.Range("A127:T129").EntireRow.Insert()

# Columns that stay constant for every "Piña / Terminal Hortofrutícola Agro Chillán" record.
$constCols = @{
    "A" = 7
    "B" = "Terminal Hortofrutícola Agro Chillán"
    "C" = "Ñuble"
    "E" = 16
    "F" = "Fruta"
    "G" = 100108
    "H" = "Tropicales y subtropicales"
    "I" = 100108005
    "J" = "Piña"
    "K" = "Caramelo"
    "R" = "Ecuador"
}

foreach ($row in 127..129) {
    foreach ($col in $constCols.Keys) {
        $ws.Range("$col$row").Value = $constCols[$col]
    }
}

# Row-specific data for the new weekly entries (fecha 44476).
$ws.Range("D127").Value = 44476
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 60
$ws.Range("N127").Value = 18000
$ws.Range("O127").Value = 19000
$ws.Range("P127").Value = 18500
$ws.Range("Q127").Value = "$/caja 12 unidades"
$ws.Range("S127").Value = 1542
$ws.Range("T127").Value = 12

$ws.Range("D128").Value = 44476
$ws.Range("L128").Value = "Segunda"
$ws.Range("M128").Value = 60
$ws.Range("N128").Value = 18000
$ws.Range("O128").Value = 19000
$ws.Range("P128").Value = 18500
$ws.Range("Q128").Value = "$/caja 14 unidades"
$ws.Range("S128").Value = 1321
$ws.Range("T128").Value = 14

$ws.Range("D129").Value = 44476
$ws.Range("L129").Value = "Tercera"
$ws.Range("M129").Value = 60
$ws.Range("N129").Value = 18000
$ws.Range("O129").Value = 19000
$ws.Range("P129").Value = 18500
$ws.Range("Q129").Value = "$/caja 16 unidades"
$ws.Range("S129").Value = 1156
$ws.Range("T129").Value = 16
